$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at position 872 (pushing the existing
# rows 872-931 down to 874-933, and growing the used range from
# A1:T931 to A1:T933).
$ws.Rows.Item(872).Insert()
$ws.Rows.Item(872).Insert()

# New row 872: Packham's Triumph / Primera, week of 2023-03-28 (45013)
$ws.Cells.Item(872, 1).Value = 8
$ws.Cells.Item(872, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(872, 3).Value = "Coquimbo"
$ws.Cells.Item(872, 4).Value = 45013
$ws.Cells.Item(872, 5).Value = 4
$ws.Cells.Item(872, 6).Value = "Fruta"
$ws.Cells.Item(872, 7).Value = 100104
$ws.Cells.Item(872, 8).Value = "Frutos de pepita"
$ws.Cells.Item(872, 9).Value = 100104005
$ws.Cells.Item(872, 10).Value = "Pera"
$ws.Cells.Item(872, 11).Value = "Packham's Triumph"
$ws.Cells.Item(872, 12).Value = "Primera"
$ws.Cells.Item(872, 13).Value = 18
$ws.Cells.Item(872, 14).Value = 220000
$ws.Cells.Item(872, 15).Value = 230000
$ws.Cells.Item(872, 16).Value = 225000
$ws.Cells.Item(872, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(872, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(872, 19).Value = 500
$ws.Cells.Item(872, 20).Value = 450

# New row 873: Packham's Triumph / Segunda, week of 2023-03-28 (45013)
$ws.Cells.Item(873, 1).Value = 8
$ws.Cells.Item(873, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(873, 3).Value = "Coquimbo"
$ws.Cells.Item(873, 4).Value = 45013
$ws.Cells.Item(873, 5).Value = 4
$ws.Cells.Item(873, 6).Value = "Fruta"
$ws.Cells.Item(873, 7).Value = 100104
$ws.Cells.Item(873, 8).Value = "Frutos de pepita"
$ws.Cells.Item(873, 9).Value = 100104005
$ws.Cells.Item(873, 10).Value = "Pera"
$ws.Cells.Item(873, 11).Value = "Packham's Triumph"
$ws.Cells.Item(873, 12).Value = "Segunda"
$ws.Cells.Item(873, 13).Value = 10
$ws.Cells.Item(873, 14).Value = 190000
$ws.Cells.Item(873, 15).Value = 200000
$ws.Cells.Item(873, 16).Value = 195000
$ws.Cells.Item(873, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(873, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(873, 19).Value = 433
$ws.Cells.Item(873, 20).Value = 450
